$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force text storage, assign, then clear the temporary format so no residual
# style/format change is left behind on the cell.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.42"
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4876"
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2876"
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06654"
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "16.80"
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07216"
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "88.68"
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.000"
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6623"
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000007813"
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.96"
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.731"
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "186.81"
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.052"
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.271"
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "155.78"
$c.ClearFormats()

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.32"
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.831"
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.401"
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.252"
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.09020"
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.930"
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05196"
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7333"
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.078"
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.693"
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01808"
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.650"
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.9210"
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.036"
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.4299"
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "104.29"
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9963"
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.716"
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1341"
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "7.269"
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.05814"
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.642"
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.416"
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.3922"
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "33.18"
$c.ClearFormats()

# Remaining cells: new text is not number-like (already has %, multiple dots,
# etc.) so a plain assignment keeps it stored as text.
$ws.Range("D2").Value = "30.628.25"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.883.85"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -4.19%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").Value = "1.874.94"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "30.556.73"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").Value = "2.120.52"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("E23").Value = "  +7.06%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +2.70%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  -6.21%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -5.76%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  -5.69%  "
$ws.Range("E38").Value = "  -3.35%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("E40").Value = "  -7.71%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -4.21%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("E46").Value = "  -7.54%  "
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -0.57%  "
